$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.046662139649411
$ws.Range("D2").Value = 1.051192063066547
$ws.Range("E2").Value = 1.044152781200023
$ws.Range("F2").Value = 1.060051279460276
$ws.Range("I2").Value = 1.037472086088742
$ws.Range("J2").Value = 1.051715340503806
$ws.Range("K2").Value = 1.053944247452335
$ws.Range("L2").Value = 1.046924634100952
$ws.Range("M2").Value = 1.062779112336001
$ws.Range("N2").Value = 1.053208896920757

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.048359903892235
$ws.Range("D3").Value = 1.052537062258058
$ws.Range("E3").Value = 1.045623329406323
$ws.Range("F3").Value = 1.061635947437379
$ws.Range("I3").Value = 1.037856112857047
$ws.Range("J3").Value = 1.053057580920572
$ws.Range("K3").Value = 1.055100474332527
$ws.Range("L3").Value = 1.048204641854546
$ws.Range("M3").Value = 1.064176184561826
$ws.Range("N3").Value = 1.054553043472872

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.049456088952768
$ws.Range("D4").Value = 1.053405036583591
$ws.Range("E4").Value = 1.046572896510557
$ws.Range("F4").Value = 1.062659339017628
$ws.Range("I4").Value = 1.038102184417826
$ws.Range("J4").Value = 1.053923419167948
$ws.Range("K4").Value = 1.055845768181098
$ws.Range("L4").Value = 1.049030414232298
$ws.Range("M4").Value = 1.065077703238852
$ws.Range("N4").Value = 1.055420111309872

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.049916367622059
$ws.Range("D5").Value = 1.053769384214792
$ws.Range("E5").Value = 1.046971631473406
$ws.Range("F5").Value = 1.063089105910124
$ws.Range("I5").Value = 1.038205056967162
$ws.Range("J5").Value = 1.054286785570952
$ws.Range("K5").Value = 1.056158412925948
$ws.Range("L5").Value = 1.049376984873235
$ws.Range("M5").Value = 1.065456117584112
$ws.Range("N5").Value = 1.055783993734848

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.04999361810956
$ws.Range("D6").Value = 1.053830527866444
$ws.Range("E6").Value = 1.047038553963648
$ws.Range("F6").Value = 1.063161238587257
$ws.Range("I6").Value = 1.038222296042579
$ws.Range("J6").Value = 1.054347759613149
$ws.Range("K6").Value = 1.056210867910987
$ws.Range("L6").Value = 1.049435141625561
$ws.Range("M6").Value = 1.06551962104653
$ws.Range("N6").Value = 1.055845054367178

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.049462241392088
$ws.Range("D7").Value = 1.053409907158233
$ws.Range("E7").Value = 1.046578226232288
$ws.Range("F7").Value = 1.062665083403953
$ws.Range("I7").Value = 1.038103561264313
$ws.Range("J7").Value = 1.053928276957304
$ws.Range("K7").Value = 1.055849948403106
$ws.Range("L7").Value = 1.049035047408388
$ws.Range("M7").Value = 1.065082761911534
$ws.Range("N7").Value = 1.055424975997847

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.047236407896613
$ws.Range("D8").Value = 1.051647099772026
$ws.Range("E8").Value = 1.044650175015217
$ws.Range("F8").Value = 1.060587244800079
$ws.Range("I8").Value = 1.037602372679028
$ws.Range("J8").Value = 1.052169518021704
$ws.Range("K8").Value = 1.054335597491024
$ws.Range("L8").Value = 1.047357737423769
$ws.Range("M8").Value = 1.063251778909966
$ws.Range("N8").Value = 1.053663719422826

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.04329538993329
$ws.Range("D9").Value = 1.048522554160126
$ws.Range("E9").Value = 1.041237119167286
$ws.Range("F9").Value = 1.056910075314856
$ws.Range("I9").Value = 1.036700536101491
$ws.Range("J9").Value = 1.049049382400808
$ws.Range("K9").Value = 1.051644815685632
$ws.Range("L9").Value = 1.044382706903922
$ws.Range("M9").Value = 1.060005917295397
$ws.Range("N9").Value = 1.050539152851464

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.040654591202623
$ws.Range("D10").Value = 1.046426676885353
$ws.Range("E10").Value = 1.038950638720668
$ws.Range("F10").Value = 1.054447368374028
$ws.Range("I10").Value = 1.036086544385615
$ws.Range("J10").Value = 1.046954563741885
$ws.Range("K10").Value = 1.049835431183151
$ws.Range("L10").Value = 1.042385728030347
$ws.Range("M10").Value = 1.057828326056706
$ws.Range("N10").Value = 1.048441359309765

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.039507725572139
$ws.Range("D11").Value = 1.045515963035232
$ws.Range("E11").Value = 1.037957790643228
$ws.Range("F11").Value = 1.053378168219417
$ws.Range("I11").Value = 1.03581760386676
$ws.Range("J11").Value = 1.046043850996945
$ws.Range("K11").Value = 1.049048144143772
$ws.Range("L11").Value = 1.041517651484496
$ws.Range("M11").Value = 1.056882021396028
$ws.Range("N11").Value = 1.047529353248268

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.039081204887266
$ws.Range("D12").Value = 1.045177193632163
$ws.Range("E12").Value = 1.037588571654229
$ws.Range("F12").Value = 1.052980581210685
$ws.Range("I12").Value = 1.035717240811908
$ws.Range("J12").Value = 1.045705012355853
$ws.Range("K12").Value = 1.048755127583209
$ws.Range("L12").Value = 1.041194691297824
$ws.Range("M12").Value = 1.056529999910944
$ws.Range("N12").Value = 1.04719003341744

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.039172718973316
$ws.Range("D13").Value = 1.045249883161425
$ws.Range("E13").Value = 1.037667790083482
$ws.Range("F13").Value = 1.053065885025658
$ws.Range("I13").Value = 1.035738790207993
$ws.Range("J13").Value = 1.045777719903929
$ws.Range("K13").Value = 1.048818007185583
$ws.Range("L13").Value = 1.041263990974427
$ws.Range("M13").Value = 1.056605533520973
$ws.Range("N13").Value = 1.047262844218574

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.03947248001373
$ws.Range("D14").Value = 1.045487970306225
$ws.Range("E14").Value = 1.037927279750423
$ws.Range("F14").Value = 1.053345312573607
$ws.Range("I14").Value = 1.035809317367308
$ws.Range("J14").Value = 1.046015853971016
$ws.Range("K14").Value = 1.049023935259385
$ws.Range("L14").Value = 1.04149096612608
$ws.Range("M14").Value = 1.056852933904674
$ws.Range("N14").Value = 1.047501316463351

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.039657102816983
$ws.Range("D15").Value = 1.045634598423087
$ws.Range("E15").Value = 1.038087102415806
$ws.Range("F15").Value = 1.053517418615446
$ws.Range("I15").Value = 1.035852709528272
$ws.Range("J15").Value = 1.046162501740679
$ws.Range("K15").Value = 1.049150736761179
$ws.Range("L15").Value = 1.041630744038176
$ws.Range("M15").Value = 1.057005295941672
$ws.Range("N15").Value = 1.047648172489667

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.040730632310183
$ws.Range("D16").Value = 1.046487049903735
$ws.Range("E16").Value = 1.039016471012997
$ws.Range("F16").Value = 1.05451826689606
$ws.Range("I16").Value = 1.036104327879708
$ws.Range("J16").Value = 1.047014926971839
$ws.Range("K16").Value = 1.049887599635547
$ws.Range("L16").Value = 1.042443267409667
$ws.Range("M16").Value = 1.057891056666718
$ws.Range("N16").Value = 1.04850180826243

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.041403112508311
$ws.Range("D17").Value = 1.047020909572387
$ws.Range("E17").Value = 1.039598684591496
$ws.Range("F17").Value = 1.055145305278837
$ws.Range("I17").Value = 1.036261334441165
$ws.Range("J17").Value = 1.047548647081183
$ws.Range("K17").Value = 1.050348786431298
$ws.Range("L17").Value = 1.042952031210995
$ws.Range("M17").Value = 1.058445754059332
$ws.Range("N17").Value = 1.049036286315547

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041795033541923
$ws.Range("D18").Value = 1.047331994325946
$ws.Range("E18").Value = 1.039938011858279
$ws.Range("F18").Value = 1.055510773818122
$ws.Range("I18").Value = 1.036352617006037
$ws.Range("J18").Value = 1.047859606375244
$ws.Range("K18").Value = 1.050617421833377
$ws.Range("L18").Value = 1.043248460000999
$ws.Range("M18").Value = 1.058768973078125
$ws.Range("N18").Value = 1.049347687207474

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041928613780098
$ws.Range("D19").Value = 1.047438014526765
$ws.Range("E19").Value = 1.040053668589053
$ws.Range("F19").Value = 1.055635343343441
$ws.Range("I19").Value = 1.036383691799315
$ws.Range("J19").Value = 1.04796557635338
$ws.Range("K19").Value = 1.050708957619013
$ws.Range("L19").Value = 1.043349479854646
$ws.Range("M19").Value = 1.058879127377131
$ws.Range("N19").Value = 1.049453807675134

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.041330995505901
$ws.Range("D20").Value = 1.046963663272077
$ws.Range("E20").Value = 1.039536246400449
$ws.Range("F20").Value = 1.055078058258167
$ws.Range("I20").Value = 1.036244519842559
$ws.Range("J20").Value = 1.047491420318106
$ws.Range("K20").Value = 1.050299343494283
$ws.Range("L20").Value = 1.042897479319617
$ws.Range("M20").Value = 1.058386274145684
$ws.Range("N20").Value = 1.048978978283903

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.039384222441862
$ws.Range("D21").Value = 1.045417873185691
$ws.Range("E21").Value = 1.037850878508414
$ws.Range("F21").Value = 1.053263040334309
$ws.Range("I21").Value = 1.035788561789782
$ws.Range("J21").Value = 1.045945744963322
$ws.Range("K21").Value = 1.048963310755166
$ws.Range("L21").Value = 1.041424141989821
$ws.Range("M21").Value = 1.056780095159179
$ws.Range("N21").Value = 1.047431107892823

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.038157171008727
$ws.Range("D22").Value = 1.044443134548553
$ws.Range("E22").Value = 1.036788720966739
$ws.Range("F22").Value = 1.052119323155166
$ws.Range("I22").Value = 1.035499181346791
$ws.Range("J22").Value = 1.044970675023039
$ws.Range("K22").Value = 1.048119915049697
$ws.Range("L22").Value = 1.040494794165824
$ws.Range("M22").Value = 1.055767201673
$ws.Range("N22").Value = 1.046454653241383

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03880794738942
$ws.Range("D23").Value = 1.044960134821237
$ws.Range("E23").Value = 1.037352031872592
$ws.Range("F23").Value = 1.052725874735739
$ws.Range("I23").Value = 1.035652844857958
$ws.Range("J23").Value = 1.045487889485325
$ws.Range("K23").Value = 1.048567338868339
$ws.Range("L23").Value = 1.040987747430226
$ws.Range("M23").Value = 1.056304446448441
$ws.Range("N23").Value = 1.046972602207531

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041363583080064
$ws.Range("D24").Value = 1.046989531358319
$ws.Range("E24").Value = 1.039564460359655
$ws.Range("F24").Value = 1.055108445133715
$ws.Range("I24").Value = 1.036252118552709
$ws.Range("J24").Value = 1.047517279712388
$ws.Range("K24").Value = 1.050321685763439
$ws.Range("L24").Value = 1.042922129973919
$ws.Range("M24").Value = 1.058413151567505
$ws.Range("N24").Value = 1.049004874401491

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.044316549440109
$ws.Range("D25").Value = 1.049332545998814
$ws.Range("E25").Value = 1.042121387351154
$ws.Range("F25").Value = 1.057862647036301
$ws.Range("I25").Value = 1.036935917328907
$ws.Range("J25").Value = 1.049858562072125
$ws.Range("K25").Value = 1.052343144958264
$ws.Range("L25").Value = 1.04515418193866
$ws.Range("M25").Value = 1.060847416592365
$ws.Range("N25").Value = 1.051349481650728
